# Daily attendance processing - reorder "Recorded By" (column G) values.
# Whenever the comma-separated list of recorders starts with "System", the
# whole list is reversed so "System" ends up last, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $reversed = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newVal = $reversed -join ", "
            $cell.Value = $newVal
        }
    }
}
